# AB#126 Sepsani funkcnich a nefunkcnich pozadavku pro diskuze s tymem
#
# The requirement "Vypadek emailove infrastruktury nezastavuje zakladni
# provoz." (a bullet under "R - Reliability (Spolehlivost)") is removed as
# a tracked deletion (paragraph mark + run text) by Ales Okon.

$d = $word.ActiveDocument

$word.UserName = "Aleš Okon"
$d.TrackRevisions = $true

$target = "Výpadek emailové infrastruktury nezastavuje základní provoz."

$rng = $d.Content
$found = $rng.Find.Execute($target, $true, $false, $false, $false, $false, `
                            $true, 1, $false, "", 0)

if ($found) {
    # $rng now spans exactly the matched text; grab its enclosing
    # paragraph and delete the whole thing (text + paragraph mark) so the
    # deletion is tracked as a removed paragraph, same as Word does when
    # you select through the pilcrow and press Delete.
    $para = $rng.Paragraphs(1)
    $para.Range.Delete()
} else {
    throw "Could not find the target sentence to delete."
}
